$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 11596.53993302689
$ws.Range("D5").Value = 11596.53993302689
$ws.Range("D9").Value = 8735.440285754541
$ws.Range("D10").Value = 8735.440285754541
$ws.Range("D14").Value = 8621.060066972881
$ws.Range("D15").Value = 8621.060066972881
